$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Powerups checklist items fixed (Speed/Size/Time Related powerups now completed) ---
$ws.Range("J26").Value = $true
$ws.Range("J27").Value = $true
$ws.Range("J28").Value = $true

# --- UI started: rename placeholder "Come up with your own" task to its real name ---
$ws.Range("C37").Value = "LERPING Or Animating Bumpers"

# --- Reflect where the author was working in the sheet when they saved ---
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("G36").Select()
